# Fruta / hortaliza, semanal
# Insert a new weekly record as row 3, pushing the existing rows 3-11 down to 4-12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 3 (shifts rows 3..11 down to 4..12)
$ws.Rows("3:3").Insert()

# Populate the newly inserted row 3 with the new record
$ws.Range("A3").Value = 12
$ws.Range("B3").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C3").Value = "Metropolitana"
$ws.Range("D3").Value = 44453
$ws.Range("E3").Value = 13
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 55
$ws.Range("K3").Value = 14000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 14455
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia de Limarí"
$ws.Range("P3").Value = 578
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
